{"js": "// Replace the multiplication-problem text in each table cell with the\n// new operands, per the commit's regenerated output. Each \"old\" string\n// appears exactly once in the document, so a plain text search + full\n// replace for each pair is sufficient and unambiguous.\nconst replacements = [\n  [\"53\u00d767=\", \"18\u00d745=\"],\n  [\"17\u00d794=\", \"32\u00d713=\"],\n  [\"25\u00d730=\", \"15\u00d766=\"],\n  [\"59\u00d785=\", \"83\u00d719=\"],\n  [\"69\u00d728=\", \"81\u00d793=\"],\n  [\"13\u00d742=\", \"25\u00d743=\"],\n  [\"19\u00d724=\", \"31\u00d768=\"],\n  [\"89\u00d767=\", \"42\u00d713=\"],\n  [\"95\u00d750=\", \"39\u00d728=\"],\n  [\"72\u00d719=\", \"98\u00d714=\"],\n  [\"14\u00d729=\", \"21\u00d786=\"],\n  [\"22\u00d747=\", \"78\u00d784=\"],\n  [\"46\u00d771=\", \"69\u00d782=\"],\n  [\"49\u00d778=\", \"85\u00d767=\"],\n  [\"93\u00d764=\", \"61\u00d760=\"],\n  [\"20\u00d787=\", \"98\u00d793=\"],\n  [\"53\u00d730=\", \"78\u00d717=\"],\n  [\"16\u00d777=\", \"44\u00d714=\"],\n  [\"91\u00d760=\", \"62\u00d795=\"],\n  [\"30\u00d776=\", \"24\u00d797=\"],\n  [\"11\u00d748=\", \"64\u00d718=\"],\n  [\"11\u00d757=\", \"47\u00d727=\"],\n  [\"78\u00d730=\", \"60\u00d781=\"],\n  [\"71\u00d762=\", \"23\u00d773=\"],\n  [\"68\u00d758=\", \"87\u00d773=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the\n# new operands, per the commit's regenerated output. Each \"old\" string\n# appears exactly once in the document, so a Find/Replace pass for each\n# pair is sufficient and unambiguous. The multiplication sign is the\n# Unicode character U+00D7 (x is the ASCII look-alike, NOT what's used\n# in the document), built via string interpolation so it round-trips as\n# a genuine \"x\" character rather than being coerced numerically by \"+\".\n$d = $word.ActiveDocument\n$x = [char]0x00D7\n\n$pairs = @(\n  @(\"53$x\" + \"67=\", \"18$x\" + \"45=\"),\n  @(\"17$x\" + \"94=\", \"32$x\" + \"13=\"),\n  @(\"25$x\" + \"30=\", \"15$x\" + \"66=\"),\n  @(\"59$x\" + \"85=\", \"83$x\" + \"19=\"),\n  @(\"69$x\" + \"28=\", \"81$x\" + \"93=\"),\n  @(\"13$x\" + \"42=\", \"25$x\" + \"43=\"),\n  @(\"19$x\" + \"24=\", \"31$x\" + \"68=\"),\n  @(\"89$x\" + \"67=\", \"42$x\" + \"13=\"),\n  @(\"95$x\" + \"50=\", \"39$x\" + \"28=\"),\n  @(\"72$x\" + \"19=\", \"98$x\" + \"14=\"),\n  @(\"14$x\" + \"29=\", \"21$x\" + \"86=\"),\n  @(\"22$x\" + \"47=\", \"78$x\" + \"84=\"),\n  @(\"46$x\" + \"71=\", \"69$x\" + \"82=\"),\n  @(\"49$x\" + \"78=\", \"85$x\" + \"67=\"),\n  @(\"93$x\" + \"64=\", \"61$x\" + \"60=\"),\n  @(\"20$x\" + \"87=\", \"98$x\" + \"93=\"),\n  @(\"53$x\" + \"30=\", \"78$x\" + \"17=\"),\n  @(\"16$x\" + \"77=\", \"44$x\" + \"14=\"),\n  @(\"91$x\" + \"60=\", \"62$x\" + \"95=\"),\n  @(\"30$x\" + \"76=\", \"24$x\" + \"97=\"),\n  @(\"11$x\" + \"48=\", \"64$x\" + \"18=\"),\n  @(\"11$x\" + \"57=\", \"47$x\" + \"27=\"),\n  @(\"78$x\" + \"30=\", \"60$x\" + \"81=\"),\n  @(\"71$x\" + \"62=\", \"23$x\" + \"73=\"),\n  @(\"68$x\" + \"58=\", \"87$x\" + \"73=\")\n)\n\nforeach ($p in $pairs) {\n  $old = $p[0]\n  $new = $p[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
